{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n\n// 1) Expand \"none of our trained models could pass it.\" into the longer\n//    sentence that lists out the classifier types that were tried.\nconst search1 = context.document.body.search(\"none of our trained models could pass it.\", { matchCase: true });\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\n    \"none of the trained classification models such as Random Forest, SVM, Na\u00efve Bayes, Logistic Regression could out perform it.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) Insert \"the \" before \"current data set\" in the closing sentence.\nconst search2 = context.document.body.search(\"I believe with current data set\", { matchCase: true });\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\n    \"I believe with the current data set\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3) Add a new, empty paragraph at the very end of the document body.\ncontext.document.body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $app / $doc all resolve; the live document is $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# 1) Expand \"none of our trained models could pass it.\" into the longer\n#    sentence that lists out the classifier types that were tried.\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.Execute(\n    \"none of our trained models could pass it.\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"none of the trained classification models such as Random Forest, SVM, Na\u00efve Bayes, Logistic Regression could out perform it.\",\n    2\n)\n\n# 2) Insert \"the \" before \"current data set\" in the closing sentence.\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.Execute(\n    \"I believe with current data set\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"I believe with the current data set\",\n    2\n)\n\n# 3) Add a new, empty paragraph at the very end of the document body.\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n$d.Save()\n"}
